# "fixed for students who have empty final marks"
# Adds a new data row (row 16) to Sheet1, cloned from row 15's layout/styles,
# for a student record whose Final Mark (column Q) is empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-establish the label merges used by every data row *before* cloning
# formats onto row 16, so the subsequent format paste (below) is what
# determines the final per-cell style - merging alone tends to stamp a
# fresh style on its anchor cell.
$ws.Range("A16:C16").Merge() | Out-Null
$ws.Range("F16:G16").Merge() | Out-Null
$ws.Range("I16:K16").Merge() | Out-Null

# Clone row 15 (formats first, then values) onto the new row 16 so every
# cell - including the blank spacer cells - keeps the same style indices.
$ws.Range("A15:Y15").Copy()
$ws.Range("A16:Y16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A15:Y15").Copy()
$ws.Range("A16:Y16").PasteSpecial(-4163)   # xlPasteValues

# New student record's encrypted student number (brand-new shared string).
$ws.Range("A16").Value = "00B197BA7753B1F2CFD57570245D6210"

# This student has no Final Mark recorded - leave Q16 blank.
$ws.Range("Q16").ClearContents()

# Widen the data columns slightly and move the active selection, matching
# the state the workbook was left in after the edit.
$ws.Columns.ColumnWidth = 17.5
$ws.Range("Q17").Select() | Out-Null
